# Insert a new event row (2024-05-03, TCT 国风动漫游戏嘉年华) as row 3 in the
# two sheets that hold event data ("展览" and "全部类型"), pushing the existing
# "HP 国风动漫游戏嘉年华" row down to row 4 and bumping its sequence number
# from 2 to 3.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Shift the existing row 3 down to row 4 and open up a blank row 3.
    $ws.Rows.Item(3).Insert()

    # New row 4 (the former row 3, "HP" event) — bump its sequence number.
    $ws.Cells.Item(4, 1).Value = 3

    # Give the new row 3's sequence-number cell (column A) the same
    # bold / centered / bordered style used by every other row's column A,
    # by copying the formatting straight from the cell below (its original
    # style, now sitting on row 4) instead of minting a near-duplicate style.
    $ws.Cells.Item(4, 1).Copy()
    $ws.Cells.Item(3, 1).PasteSpecial(-4122)  # xlPasteFormats

    # New row 3: the TCT event.
    $ws.Cells.Item(3, 1).Value = 2

    # Assign the date as literal text (not an auto-converted date serial):
    # a leading apostrophe forces Excel to store it verbatim as a string,
    # then resetting the style back to Normal drops the quote-prefix flag
    # so the cell ends up with the same plain (unstyled) look as its peers.
    $ws.Cells.Item(3, 2).Value = "'2024-05-03"
    $ws.Cells.Item(3, 2).Style = "Normal"

    $ws.Cells.Item(3, 3).Value = "丽水·首届TCT国风动漫游戏嘉年华"
    $ws.Cells.Item(3, 4).Value = "括苍路493号油泵厂山顶通用设备厂区块3号楼 中国国际摄影节展览馆"
    $ws.Cells.Item(3, 5).Value = "2024.05.03 10:00-05.04 17:00"
    $ws.Cells.Item(3, 6).Value = 101
    $ws.Cells.Item(3, 7).Value = 29.9
    $ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84156"
    $ws.Cells.Item(3, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/rpRBCHaU1712892375435.jpeg"
}
